$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("航天发展", "襄阳轴承", "峨眉山A"),
    @("雷科防务", "峨眉山Ａ", "航天发展"),
    @("乾照光电", "航天发展", "襄阳轴承"),
    @("平潭发展", "雷科防务", "雷科防务"),
    @("国晟科技", "乾照光电", "平潭发展"),
    @("通宇通讯", "华电新能", "东方精工"),
    @("襄阳轴承", "平潭发展", "国晟科技"),
    @("天风证券", "通宇通讯", "实达集团"),
    @("实达集团", "国晟科技", "华电新能"),
    @("华电新能", "天风证券", "乾照光电"),
    @("海欣食品", "海欣食品", "通宇通讯"),
    @("航天动力", "中科曙光", "航天动力"),
    @("赛微电子", "上汽集团", "合富中国"),
    @("银河电子", "航天动力", "天风证券"),
    @("海王生物", "银河电子", "银河电子"),
    @("航天环宇", "实达集团", "海欣食品"),
    @("中科曙光", "东方精工", "海王生物"),
    @("中际旭创", "顺灏股份", "中科曙光"),
    @("光库科技", "海王生物", "特发信息"),
    @("峨眉山Ａ", "北方稀土", "多氟多")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}
